# Updated symbol list on Wed Dec 28 05:22:58 UTC 2022 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "244.62"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "23.82"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.332"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.05784"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "6.483"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.337"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.8108"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.8877"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1392"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07361"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03090"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03059"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09355"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.870"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.001542"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.04719"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0006006"
$ws.Range("E18").Value = "17OneONE"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.006099"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.001294"
$ws.Range("B21").Value = "NitroEx"
$ws.Range("C21").Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.00008810"
$ws.Range("E21").Value = "20NitroExNTXBestin24h"
$ws.Range("B22").Value = "LEO"
$ws.Range("C22").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.581"
$ws.Range("E22").Value = "21LEOLEO"
$ws.Range("B23").Value = "BTSEToken"
$ws.Range("C23").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.143"
$ws.Range("E23").Value = "22BTSETokenBTSE"
$ws.Range("B24").Value = "BitpandaEcosystemToken"
$ws.Range("C24").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.3180"
$ws.Range("E24").Value = "23BitpandaEcosystemTokenBEST"
$ws.Range("B25").Value = "ProBitToken"
$ws.Range("C25").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1318"
$ws.Range("E25").Value = "24ProBitTokenPROB"
$ws.Range("B26").Value = "AAXToken"
$ws.Range("C26").Value = "https://coinranking.com/coin/LNePqkIhk+aaxtoken-aab"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.2000"
$ws.Range("E26").Value = "25AAXTokenAAB"
$ws.Range("B27").Value = "HotbitToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.004618"
$ws.Range("E27").Value = "26HotbitTokenHTB"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.03807"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006367"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1055"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002753"
$ws.Range("E43").Value = "42CEJICEJIWorstin24h"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.007611"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005448"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00000000751"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5506"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.001841"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00002102"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0002002"
